$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -eq $val) { continue }

    # Only rewrite "Recorded By" entries that start with a leading "System, "
    # token and do not involve the admin@admin.com account.
    if ($val.StartsWith("System, ") -and -not $val.Contains("admin@admin.com")) {
        $parts = $val -split ", "
        $rest = $parts[1..($parts.Count - 1)]
        $newVal = ($rest + $parts[0]) -join ", "
        $cell.Value2 = $newVal
    }
}
